$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 9: copy formatting (date style) from A8 down to A9, then set values
$ws.Range("A8").Copy($ws.Range("A9"))
$ws.Range("A9").Value = 42149

$ws.Range("B9").Value = "Comment for the demo project: Create, update, delete and get comment on service"

$ws.Range("B9").Select()
